$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting (bold + border + alignment) from an existing
# header cell onto the three new header cells, then set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every player row.
for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 30).Value = 82   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 80   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
